# Auto-generated edit script applying the Tonberry_Profits.xlsx diff
# Updates LeveProfit-related computed columns (H-N) on several rows
# across all 8 job sheets. All target cells hold plain numeric values
# (no formulas in this workbook), so we just overwrite/clear values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M12").Value = -962
$ws.Range("I12").Value = 1132
$ws.Range("K12").Value = 1132
$ws.Range("H12").Value = 1235
$ws.Range("J17").Value = 3379.3
$ws.Range("L17").Value = 10137.9
$ws.Range("H17").Value = 3850.5454
$ws.Range("N17").Value = -10473.9
$ws.Range("I43").Value = 2576.6667
$ws.Range("K43").Value = 2576.6667
$ws.Range("H43").Value = 2038.3334
$ws.Range("M43").Value = -2507.6667
$ws.Range("L70").Value = 3300
$ws.Range("N70").Value = -3840
$ws.Range("J70").Value = 1100
$ws.Range("H70").Value = 1100
$ws.Range("L73").Value = 3300
$ws.Range("H73").Value = 1100
$ws.Range("N73").Value = -5172
$ws.Range("J73").Value = 1100
$ws.Range("M76").Value = -5028.75
$ws.Range("H76").Value = 5323
$ws.Range("N76").Value = -5925.3335
$ws.Range("K76").Value = 5343.75
$ws.Range("J76").Value = 5295.3335
$ws.Range("I76").Value = 5343.75
$ws.Range("L76").Value = 5295.3335
$ws.Range("M79").Value = -4251.75
$ws.Range("I79").Value = 5343.75
$ws.Range("K79").Value = 5343.75
$ws.Range("H79").Value = 5323
$ws.Range("N79").Value = -7479.3335
$ws.Range("J79").Value = 5295.3335
$ws.Range("L79").Value = 5295.3335
$ws.Range("H86").Value = 1639.8
$ws.Range("I86").Value = 1599.7778
$ws.Range("K86").Value = 1599.7778
$ws.Range("M86").Value = -476.7778000000001
$ws.Range("I89").Value = 1599.7778
$ws.Range("K89").Value = 7998.889
$ws.Range("H89").Value = 1639.8
$ws.Range("M89").Value = -2382.889
$ws.Range("I106").Value = 2311.1
$ws.Range("K106").Value = 2311.1
$ws.Range("H106").Value = 3105.52
$ws.Range("M106").Value = -1680.1
$ws.Range("L129").Value = 4914.2856
$ws.Range("N129").Value = -14914.2856
$ws.Range("H129").Value = 1581.7273
$ws.Range("J129").Value = 1638.0952
$ws.Range("K132").Value = 3225.3915
$ws.Range("M132").Value = -695.3914999999997
$ws.Range("I132").Value = 1075.1305
$ws.Range("H132").Value = 1154.3265
$ws.Range("K137").Value = 4063.3638
$ws.Range("M137").Value = -1513.3638
$ws.Range("I137").Value = 1354.4546
$ws.Range("H137").Value = 43210.457
$ws.Range("J138").Value = 3462.5098
$ws.Range("N138").Value = -20667.5294
$ws.Range("L138").Value = 10387.5294
$ws.Range("H138").Value = 3800.493

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17458.164
$ws.Range("J32").Value = 25824.262
$ws.Range("L32").Value = 25824.262
$ws.Range("M32").Value = -12797.978
$ws.Range("N32").Value = -26398.262
$ws.Range("I32").Value = 13084.978
$ws.Range("K32").Value = 13084.978
$ws.Range("M45").Value = -3855.6665
$ws.Range("K45").Value = 4232.6665
$ws.Range("H45").Value = 2516.1667
$ws.Range("I45").Value = 4232.6665
$ws.Range("H61").Value = 25591.922
$ws.Range("M61").Value = -32834.96
$ws.Range("I61").Value = 33046.96
$ws.Range("K61").Value = 33046.96
$ws.Range("H74").Value = 1404.5526
$ws.Range("M74").Value = 57.92589999999996
$ws.Range("K74").Value = 816.0741
$ws.Range("I74").Value = 816.0741
$ws.Range("M77").Value = 287.6295
$ws.Range("H77").Value = 1404.5526
$ws.Range("I77").Value = 816.0741
$ws.Range("K77").Value = 4080.3705
$ws.Range("L88").Value = 4566.6665
$ws.Range("N88").Value = -5378.6665
$ws.Range("H88").Value = 4288.25
$ws.Range("J88").Value = 4566.6665
$ws.Range("H91").Value = 4288.25
$ws.Range("L91").Value = 4566.6665
$ws.Range("J91").Value = 4566.6665
$ws.Range("N91").Value = -7374.6665
$ws.Range("K132").Value = 4289.5587
$ws.Range("M132").Value = -1759.5587
$ws.Range("I132").Value = 1429.8529
$ws.Range("H132").Value = 1827.8529
$ws.Range("I136").Value = 33046.96
$ws.Range("M136").Value = -96590.88
$ws.Range("K136").Value = 99140.88
$ws.Range("H136").Value = 25591.922

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("I10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H64").Value = 184.55556
$ws.Range("L64").Value = 88.25
$ws.Range("N64").Value = -538.25
$ws.Range("J64").Value = 88.25
$ws.Range("J67").Value = 88.25
$ws.Range("N67").Value = -1648.25
$ws.Range("H67").Value = 184.55556
$ws.Range("L67").Value = 88.25
$ws.Range("H86").Value = 1002003.5
$ws.Range("N86").Value = -1004249.5
$ws.Range("J86").Value = 1002003.5
$ws.Range("L86").Value = 1002003.5
$ws.Range("N89").Value = -5021249.5
$ws.Range("H89").Value = 1002003.5
$ws.Range("L89").Value = 5010017.5
$ws.Range("J89").Value = 1002003.5
$ws.Range("N134").Value = -16509.75
$ws.Range("H134").Value = 4572.3477
$ws.Range("J134").Value = 3813.25
$ws.Range("L134").Value = 11439.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J31").Value = 5482.625
$ws.Range("K31").Value = 1613.65
$ws.Range("I31").Value = 1613.65
$ws.Range("M31").Value = -1318.65
$ws.Range("L31").Value = 5482.625
$ws.Range("H31").Value = 3333.1943
$ws.Range("N31").Value = -6072.625
$ws.Range("N34").Value = -5886.625
$ws.Range("J34").Value = 5482.625
$ws.Range("L34").Value = 5482.625
$ws.Range("I34").Value = 1613.65
$ws.Range("K34").Value = 1613.65
$ws.Range("H34").Value = 3333.1943
$ws.Range("M34").Value = -1411.65
$ws.Range("H86").Value = 9597.429
$ws.Range("I86").Value = 2366.4
$ws.Range("K86").Value = 2366.4
$ws.Range("N86").Value = -15860.667
$ws.Range("J86").Value = 13614.667
$ws.Range("M86").Value = -1243.4
$ws.Range("L86").Value = 13614.667
$ws.Range("I89").Value = 2366.4
$ws.Range("N89").Value = -79305.33499999999
$ws.Range("K89").Value = 11832
$ws.Range("H89").Value = 9597.429
$ws.Range("L89").Value = 68073.33499999999
$ws.Range("J89").Value = 13614.667
$ws.Range("M89").Value = -6216
$ws.Range("M134").Value = -639.0935999999997
$ws.Range("H134").Value = 1212.3096
$ws.Range("I134").Value = 1058.0312
$ws.Range("K134").Value = 3174.0936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2262.0908
$ws.Range("K39").Value = 5989.799999999999
$ws.Range("M39").Value = -5695.799999999999
$ws.Range("I39").Value = 1996.6
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("H113").Value = 26383.932
$ws.Range("J113").Value = 1491.7354
$ws.Range("N113").Value = -8815.206200000001
$ws.Range("L113").Value = 4475.206200000001
$ws.Range("L122").Value = 12079.8
$ws.Range("I122").Value = 740
$ws.Range("M122").Value = -4210
$ws.Range("N122").Value = -16979.8
$ws.Range("H122").Value = 1191.65
$ws.Range("J122").Value = 1342.2
$ws.Range("K122").Value = 6660
$ws.Range("H131").Value = 50648.125
$ws.Range("J131").Value = 50648.125
$ws.Range("L131").Value = 151944.375
$ws.Range("N131").Value = -162024.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L46").Value = 28666.666
$ws.Range("J46").Value = 28666.666
$ws.Range("H46").Value = 28666.666
$ws.Range("N46").Value = -28978.666
$ws.Range("I80").Value = 2424.6667
$ws.Range("H80").Value = 2454.9473
$ws.Range("K80").Value = 2424.6667
$ws.Range("M80").Value = -1426.6667
$ws.Range("I83").Value = 2424.6667
$ws.Range("K83").Value = 12123.3335
$ws.Range("H83").Value = 2454.9473
$ws.Range("M83").Value = -7131.333500000001
$ws.Range("M97").Value = -837
$ws.Range("H97").Value = 1333
$ws.Range("I97").Value = 1333
$ws.Range("K97").Value = 1333
$ws.Range("N132").Value = -15499.667
$ws.Range("L132").Value = 10439.667
$ws.Range("K132").Value = 2696714.7
$ws.Range("M132").Value = -2694184.7
$ws.Range("I132").Value = 898904.9
$ws.Range("H132").Value = 634681.1
$ws.Range("J132").Value = 3479.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 1976.3125
$ws.Range("I7").Value = 1976.3125
$ws.Range("M7").Value = -1864.3125
$ws.Range("H7").Value = 2428.1482
$ws.Range("N34").Value = -60368
$ws.Range("J34").Value = 60024
$ws.Range("L34").Value = 60024
$ws.Range("H34").Value = 60024
$ws.Range("I126").Value = 1976.3125
$ws.Range("H126").Value = 2428.1482
$ws.Range("M126").Value = -3458.9375
$ws.Range("K126").Value = 5928.9375
$ws.Range("N132").Value = -34690.667
$ws.Range("L132").Value = 29630.667
$ws.Range("K132").Value = 14726.4
$ws.Range("M132").Value = -12196.4
$ws.Range("I132").Value = 4908.8
$ws.Range("H132").Value = 6771.8335
$ws.Range("J132").Value = 9876.888999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 12053.714
$ws.Range("L96").Value = 16035.2
$ws.Range("N96").Value = -18781.2
$ws.Range("J96").Value = 16035.2
$ws.Range("K100").Value = 1042
$ws.Range("I100").Value = 521
$ws.Range("M100").Value = -501
$ws.Range("L100").Value = 1627.2
$ws.Range("J100").Value = 813.6
$ws.Range("H100").Value = 602.2778
$ws.Range("N100").Value = -2709.2
$ws.Range("I126").Value = 2857.7896
$ws.Range("H126").Value = 2708.4285
$ws.Range("M126").Value = -6103.3688
$ws.Range("K126").Value = 8573.3688
